# Apply refreshed market data from the Tue Dec 12 04:26:04 UTC 2023 GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds "Price" as plain text (values like "41.717.02" use "." as a thousands
# separator so Excel cannot read them as numbers). A leading apostrophe forces the new
# price to be stored as text too, matching the original formatting; Excel strips the
# apostrophe quote-prefix marker from the value actually stored in the cell.
function Set-PriceText($row, $value) {
    $ws.Cells.Item($row, 4).Value = "'$value"
}

Set-PriceText 2 "41.730.85"
$ws.Cells.Item(2, 5).Value = "  -1.15%  "
Set-PriceText 3 "2.231.47"
$ws.Cells.Item(3, 5).Value = "  -0.60%  "
$ws.Cells.Item(4, 5).Value = "  -0.17%  "
Set-PriceText 5 "251.71"
$ws.Cells.Item(5, 5).Value = "  +8.24%  "
Set-PriceText 6 "0.630"
$ws.Cells.Item(6, 5).Value = "  -0.96%  "
Set-PriceText 7 "71.19"
$ws.Cells.Item(7, 5).Value = "  +0.72%  "
$ws.Cells.Item(8, 5).Value = "  +0.03%  "
Set-PriceText 9 "0.562"
$ws.Cells.Item(9, 5).Value = "  -0.07%  "
Set-PriceText 10 "41.83"
$ws.Cells.Item(10, 5).Value = "  +16.23%  "
Set-PriceText 11 "0.0965"
$ws.Cells.Item(11, 5).Value = "  -4.84%  "
Set-PriceText 12 "58.53"
$ws.Cells.Item(12, 5).Value = "  +0.07%  "
Set-PriceText 13 "0.105"
$ws.Cells.Item(13, 5).Value = "  +0.55%  "
Set-PriceText 14 "6.91"
$ws.Cells.Item(14, 5).Value = "  +0.53%  "
Set-PriceText 15 "2.558.20"
$ws.Cells.Item(15, 5).Value = "  -0.92%  "
Set-PriceText 16 "14.97"
$ws.Cells.Item(16, 5).Value = "  -1.06%  "
Set-PriceText 17 "0.857"
$ws.Cells.Item(17, 5).Value = "  -1.67%  "
Set-PriceText 18 "2.233.03"
$ws.Cells.Item(18, 5).Value = "  -0.63%  "
Set-PriceText 19 "41.677.03"
$ws.Cells.Item(19, 5).Value = "  -1.05%  "
Set-PriceText 20 "0.0₃0968"
$ws.Cells.Item(20, 5).Value = "  -2.96%  "
Set-PriceText 21 "72.99"
$ws.Cells.Item(21, 5).Value = "  -1.09%  "
Set-PriceText 22 "6.17"
$ws.Cells.Item(22, 5).Value = "  -0.98%  "
Set-PriceText 23 "2.25"
$ws.Cells.Item(23, 5).Value = "  +14.58%  "
Set-PriceText 24 "234.42"
$ws.Cells.Item(24, 5).Value = "  -1.96%  "
$ws.Cells.Item(25, 5).Value = "  +0.10%  "
$ws.Cells.Item(26, 5).Value = "  +2.13%  "
Set-PriceText 27 "2.52"
$ws.Cells.Item(27, 5).Value = "  +7.01%  "
Set-PriceText 28 "10.18"
$ws.Cells.Item(28, 5).Value = "  -0.14%  "
Set-PriceText 29 "2.19"
$ws.Cells.Item(29, 5).Value = "  +1.34%  "
Set-PriceText 30 "170.21"
$ws.Cells.Item(30, 5).Value = "  +1.10%  "
Set-PriceText 31 "20.64"
$ws.Cells.Item(31, 5).Value = "  -0.51%  "
$ws.Cells.Item(32, 5).Value = "  +0.09%  "
$ws.Cells.Item(33, 5).Value = "  -1.85%  "
Set-PriceText 34 "5.49"
$ws.Cells.Item(34, 5).Value = "  +2.23%  "
Set-PriceText 35 "0.0720"
$ws.Cells.Item(35, 5).Value = "  -1.01%  "
Set-PriceText 36 "26.58"
$ws.Cells.Item(36, 5).Value = "  +17.57%  "
Set-PriceText 37 "4.67"
$ws.Cells.Item(37, 5).Value = "  -2.80%  "
Set-PriceText 38 "4.09"
$ws.Cells.Item(38, 5).Value = "  +12.86%  "
Set-PriceText 39 "0.0288"
$ws.Cells.Item(39, 5).Value = "  +7.62%  "
Set-PriceText 40 "70.59"
$ws.Cells.Item(40, 5).Value = "  +4.56%  "
Set-PriceText 41 "2.28"
$ws.Cells.Item(41, 5).Value = "  +1.49%  "
Set-PriceText 42 "6.02"
$ws.Cells.Item(42, 5).Value = "  -1.80%  "
Set-PriceText 43 "0.212"
$ws.Cells.Item(43, 5).Value = "  +12.65%  "
$ws.Cells.Item(48, 5).Value = "  +0.77%  "
$ws.Cells.Item(49, 5).Value = "  -0.07%  "
$ws.Cells.Item(50, 5).Value = "  +6.22%  "
$ws.Cells.Item(51, 5).Value = "  +1.04%  "

# Rows 44-47 were re-ranked by 24h change: swap coin name, link, price and change together
$ws.Cells.Item(44, 2).Value = "FTXToken"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-PriceText 44 "5.14"
$ws.Cells.Item(44, 5).Value = "  +4.20%  "
$ws.Cells.Item(45, 2).Value = "Celestia"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-PriceText 45 "11.78"
$ws.Cells.Item(45, 5).Value = "  +13.40%  "
$ws.Cells.Item(46, 2).Value = "FraxShare"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-PriceText 46 "8.81"
$ws.Cells.Item(46, 5).Value = "  -1.35%  "
$ws.Cells.Item(47, 2).Value = "SynthetixNetwork"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
Set-PriceText 47 "4.80"
$ws.Cells.Item(47, 5).Value = "  +8.48%  "

Write-Host "Applied cryptos list update"
